$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "PRODUCTO"
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").HorizontalAlignment = -4108

$lastRow = 213
$rng = $ws.Range("K2:K" + $lastRow)
$rng.Value = "CEBADA"
